$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.676.85'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '3.465.35'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''414.67'
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").Value = '''130.24'
$ws.Range("E6").Value = '  +1.52%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("D9").Value = '''0.727'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("D11").Value = '''42.67'
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").Value = '''9.75'
$ws.Range("E12").Value = '  +7.13%  '
$ws.Range("D13").Value = '''0.0000219'
$ws.Range("E13").Value = '  +6.41%  '
$ws.Range("D14").Value = '4.013.85'
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").Value = '''0.140'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").Value = '''20.50'
$ws.Range("E16").Value = '  -3.82%  '
$ws.Range("D17").Value = '3.469.81'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").Value = '''12.76'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '62.628.69'
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").Value = '''469.99'
$ws.Range("E21").Value = '  +5.84%  '
$ws.Range("D22").Value = '''90.68'
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").Value = '''3.26'
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("D24").Value = '''13.39'
$ws.Range("E24").Value = '  +3.38%  '
$ws.Range("D25").Value = '''10.52'
$ws.Range("E25").Value = '  +20.50%  '
$ws.Range("D26").Value = '''3.30'
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").Value = '''7.61'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").Value = '''12.01'
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("E32").Value = '  -1.79%  '
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").Value = '''41.15'
$ws.Range("E34").Value = '  -3.93%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '''58.81'
$ws.Range("E36").Value = '  +8.99%  '
$ws.Range("D37").Value = '''0.0491'
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").Value = '''3.06'
$ws.Range("E39").Value = '  +4.45%  '
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  +6.63%  '
$ws.Range("D44").Value = '''145.25'
$ws.Range("E44").Value = '  +2.41%  '
$ws.Range("D45").Value = '''4.36'
$ws.Range("E45").Value = '  +2.71%  '
$ws.Range("E46").Value = '  +4.04%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").Value = '''2.40'
$ws.Range("E47").Value = '  +10.84%  '
$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").Value = '0.0₃0569'
$ws.Range("E48").Value = '  +39.66%  '
$ws.Range("D49").Value = '''16.38'
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").Value = '''22.18'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  -1.30%  '
